$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20 mirrors the existing date-formatted rows (copy format from row 19)
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -0.7200474048664085
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -0.2284091334091687
